# CDS Input file updates
# Updates the Cypher queries stored in column B (per-tab query) for the
# ParticipantsTab, SamplesTab and FilesTab rows, refreshes the shared
# "StatQuery" text in column C, resizes row heights to fit the longer
# wrapped text, widens column A, and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New query text for each tab --------------------------------------

$participantsQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE g.platform in ['Illumina HiSeq 2500']
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id limit 100
'@

$samplesQuery = @'
MATCH (samp:sample)-->(p:participant)-->(s:study)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE g.platform in ['Illumina HiSeq 2500']
WITH DISTINCT s, p, samp
RETURN
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(p.participant_id,'') as `Participant ID`,
    coalesce(s.study_name, '') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(samp.sample_tumor_status,'') as `Tumor`,
    coalesce(samp.sample_type,'') as `Analyte Type`
ORDER BY samp.sample_id limit 100
'@

$filesQuery = @'
MATCH (f:file)-->(s:study)
OPTIONAL MATCH (samp:sample)<--(f)
OPTIONAL MATCH (samp)-->(p:participant)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE g.platform in ['Illumina HiSeq 2500']
WITH DISTINCT f, s, p, samp
RETURN
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name,'') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id, '') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER BY f.file_name limit 100
'@

$statQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.platform in ['Illumina HiSeq 2500']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH DISTINCT samp,s,p,f
RETURN
    count(distinct s) AS Studies,
    count(distinct p) AS Participants,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Files`
'@

# --- Row 2: ParticipantsTab --------------------------------------------
$ws.Range("B2").Value = $participantsQuery
$ws.Range("C2").Value = $statQuery
$ws.Rows.Item(2).RowHeight = 263.5

# --- Row 3: SamplesTab ---------------------------------------------------
$ws.Range("B3").Value = $samplesQuery
$ws.Range("C3").Value = $statQuery
$ws.Rows.Item(3).RowHeight = 232.5

# --- Row 4: FilesTab ------------------------------------------------------
$ws.Range("B4").Value = $filesQuery
$ws.Range("C4").Value = $statQuery
$ws.Rows.Item(4).RowHeight = 248

# --- Column A width (no longer "best fit") --------------------------------
# Target stored width is 25.54296875 characters; the COM ColumnWidth
# property is quantized to a 1/6-character pixel grid, so feed it the
# calibrated input that lands closest to the target after conversion.
$ws.Columns.Item(1).ColumnWidth = 24.666666666666668

# --- Active cell / selection -----------------------------------------------
$ws.Range("B3").Select()

Write-Host "edit applied"
